$d = $word.ActiveDocument

# 1) "Critério" sentence: split the trailing formula onto its own line
#    by inserting a manual line break before "NF = NOTA x % FREQ."
$d.Content.Find.Execute(
    "abaixo:NF = NOTA",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "abaixo:^lNF = NOTA",
    2
)

# 2) Bibliografia paragraph: insert a blank line (two manual line breaks)
#    between each reference entry.
$d.Content.Find.Execute(
    "2011.Encyclopedia",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "2011.^l^lEncyclopedia",
    2
)

$d.Content.Find.Execute(
    "2006.Kirk,",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "2006.^l^lKirk,",
    2
)

$d.Content.Find.Execute(
    "1984.Manual",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "1984.^l^lManual",
    2
)

$d.Content.Find.Execute(
    "2007.Shreve,",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "2007.^l^lShreve,",
    2
)

$d.Content.Find.Execute(
    "c1997.Revistas:",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "c1997.^l^lRevistas:",
    2
)

$d.Content.Find.Execute(
    "Revistas:Brazilian",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Revistas:^l^lBrazilian",
    2
)
